$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 507. This shifts the existing rows 507..628
# down to 508..629, matching the target diff.
$ws.Rows.Item(507).Insert()

# Populate the newly inserted row 507 with the new record's values.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R are constant across all the
# rows in this block, so reuse them; only D, J, K, L, M, P change.
$ws.Cells.Item(507, 1).Value = 3
$ws.Cells.Item(507, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(507, 3).Value = "Coquimbo"
$ws.Cells.Item(507, 4).Value = 45204
$ws.Cells.Item(507, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(507, 5).Value = 5
$ws.Cells.Item(507, 6).Value = 100112009
$ws.Cells.Item(507, 7).Value = "Acelga"
$ws.Cells.Item(507, 8).Value = "Sin especificar"
$ws.Cells.Item(507, 9).Value = "Primera"
$ws.Cells.Item(507, 10).Value = 220
$ws.Cells.Item(507, 11).Value = 3000
$ws.Cells.Item(507, 12).Value = 3500
$ws.Cells.Item(507, 13).Value = 3250
$ws.Cells.Item(507, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(507, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(507, 16).Value = 542
$ws.Cells.Item(507, 17).Value = 6
$ws.Cells.Item(507, 18).Value = "Hortaliza"
